$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='62.166.74'; Text=$false}
    @{Cell='D3'; Value='2.440.96'; Text=$false}
    @{Cell='E3'; Value='  -1.37%  '; Text=$false}
    @{Cell='E4'; Value='  -0.12%  '; Text=$false}
    @{Cell='D5'; Value='578.89'; Text=$true}
    @{Cell='E5'; Value='  +0.51%  '; Text=$false}
    @{Cell='D6'; Value='142.86'; Text=$true}
    @{Cell='E6'; Value='  -4.08%  '; Text=$false}
    @{Cell='E7'; Value='  -0.01%  '; Text=$false}
    @{Cell='E8'; Value='  -2.56%  '; Text=$false}
    @{Cell='D9'; Value='2.438.59'; Text=$false}
    @{Cell='E9'; Value='  -1.40%  '; Text=$false}
    @{Cell='E10'; Value='  -4.95%  '; Text=$false}
    @{Cell='E11'; Value='  +1.20%  '; Text=$false}
    @{Cell='D12'; Value='5.19'; Text=$true}
    @{Cell='E12'; Value='  -2.23%  '; Text=$false}
    @{Cell='E13'; Value='  -3.70%  '; Text=$false}
    @{Cell='D14'; Value='26.39'; Text=$true}
    @{Cell='E14'; Value='  -3.29%  '; Text=$false}
    @{Cell='D15'; Value='0.0000172'; Text=$true}
    @{Cell='E15'; Value='  -5.78%  '; Text=$false}
    @{Cell='D16'; Value='2.867.99'; Text=$false}
    @{Cell='E16'; Value='  -1.89%  '; Text=$false}
    @{Cell='D17'; Value='62.205.92'; Text=$false}
    @{Cell='E17'; Value='  -2.09%  '; Text=$false}
    @{Cell='D18'; Value='2.436.96'; Text=$false}
    @{Cell='E18'; Value='  -2.12%  '; Text=$false}
    @{Cell='D19'; Value='10.94'; Text=$true}
    @{Cell='E19'; Value='  -4.68%  '; Text=$false}
    @{Cell='D20'; Value='7.09'; Text=$true}
    @{Cell='E20'; Value='  -4.64%  '; Text=$false}
    @{Cell='D21'; Value='329.19'; Text=$true}
    @{Cell='E21'; Value='  -0.78%  '; Text=$false}
    @{Cell='D22'; Value='4.11'; Text=$true}
    @{Cell='E22'; Value='  -3.05%  '; Text=$false}
    @{Cell='D23'; Value='1.95'; Text=$true}
    @{Cell='E23'; Value='  -8.05%  '; Text=$false}
    @{Cell='E24'; Value='  -0.07%  '; Text=$false}
    @{Cell='D25'; Value='65.49'; Text=$true}
    @{Cell='E25'; Value='  -0.82%  '; Text=$false}
    @{Cell='D26'; Value='9.30'; Text=$true}
    @{Cell='E26'; Value='  +0.83%  '; Text=$false}
    @{Cell='D27'; Value='627.47'; Text=$true}
    @{Cell='E27'; Value='  -0.32%  '; Text=$false}
    @{Cell='D28'; Value='2.560.34'; Text=$false}
    @{Cell='E28'; Value='  -1.72%  '; Text=$false}
    @{Cell='B29'; Value='PEPE'; Text=$false}
    @{Cell='C29'; Value='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; Text=$false}
    @{Cell='D29'; Value='0.0₃0948'; Text=$false}
    @{Cell='E29'; Value='  -9.89%  '; Text=$false}
    @{Cell='B30'; Value='Binance-PegBSC-USD'; Text=$false}
    @{Cell='C30'; Value='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; Text=$false}
    @{Cell='D30'; Value='1.00'; Text=$true}
    @{Cell='E30'; Value='  -0.11%  '; Text=$false}
    @{Cell='D31'; Value='1.43'; Text=$true}
    @{Cell='E31'; Value='  -7.79%  '; Text=$false}
    @{Cell='D32'; Value='7.99'; Text=$true}
    @{Cell='E32'; Value='  -5.22%  '; Text=$false}
    @{Cell='E33'; Value='  -1.37%  '; Text=$false}
    @{Cell='E34'; Value='  -2.02%  '; Text=$false}
    @{Cell='D35'; Value='4.93'; Text=$true}
    @{Cell='E35'; Value='  -6.38%  '; Text=$false}
    @{Cell='E36'; Value='  +0.10%  '; Text=$false}
    @{Cell='E37'; Value='  -8.37%  '; Text=$false}
    @{Cell='D38'; Value='0.375'; Text=$true}
    @{Cell='E38'; Value='  -2.45%  '; Text=$false}
    @{Cell='D39'; Value='150.06'; Text=$true}
    @{Cell='E39'; Value='  +2.71%  '; Text=$false}
    @{Cell='D40'; Value='18.28'; Text=$true}
    @{Cell='E40'; Value='  -3.44%  '; Text=$false}
    @{Cell='D41'; Value='5.22'; Text=$true}
    @{Cell='E41'; Value='  -5.97%  '; Text=$false}
    @{Cell='E42'; Value='  -3.97%  '; Text=$false}
    @{Cell='D43'; Value='42.73'; Text=$true}
    @{Cell='E43'; Value='  +1.88%  '; Text=$false}
    @{Cell='E44'; Value='  -0.01%  '; Text=$false}
    @{Cell='E45'; Value='  -10.24%  '; Text=$false}
    @{Cell='D46'; Value='142.75'; Text=$true}
    @{Cell='E46'; Value='  -5.53%  '; Text=$false}
    @{Cell='D47'; Value='3.64'; Text=$true}
    @{Cell='E47'; Value='  -4.11%  '; Text=$false}
    @{Cell='E48'; Value='  -3.72%  '; Text=$false}
    @{Cell='D49'; Value='0.598'; Text=$true}
    @{Cell='E49'; Value='  -1.61%  '; Text=$false}
    @{Cell='D50'; Value='19.54'; Text=$true}
    @{Cell='E50'; Value='  -9.58%  '; Text=$false}
    @{Cell='D51'; Value='0.0₆0231'; Text=$false}
    @{Cell='E51'; Value='  +1.86%  '; Text=$false}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Text) {
        $rng.Value = "'" + $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
